# Apply cryptos list update (diff-driven cell edits)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.628.74'
$ws.Range('E2').Value = '  -0.30%  '
$ws.Range('D3').Value = '3.323.01'
$ws.Range('E3').Value = '  -0.01%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '578.80'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.02%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '174.87'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.34%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.587'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.59%  '
$ws.Range('D9').Value = '3.320.05'
$ws.Range('E9').Value = '  +0.08%  '
$ws.Range('E10').Value = '  -0.59%  '
$ws.Range('E11').Value = '  -0.75%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '45.27'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.24%  '
$ws.Range('E13').Value = '  -1.87%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '658.83'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.88%  '
$ws.Range('D15').Value = '3.868.23'
$ws.Range('E15').Value = '  +0.30%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '8.38'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.02%  '
$ws.Range('D17').Value = '67.663.80'
$ws.Range('E17').Value = '  -0.41%  '
$ws.Range('E18').Value = '  -0.77%  '
$ws.Range('D19').Value = '3.325.15'
$ws.Range('E19').Value = '  +0.06%  '
$ws.Range('E20').Value = '  -2.02%  '
$ws.Range('E21').Value = '  +0.13%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.885'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.94%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.37'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +6.69%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '17.01'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.43%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '98.74'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.90%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.83'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.02%  '
$ws.Range('E27').Value = '  -4.07%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.23'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.75%  '
$ws.Range('E29').Value = '  +2.44%  '
$ws.Range('E30').Value = '  -2.50%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.23'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +7.84%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '569.95'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.26%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '10.92'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.16%  '
$ws.Range('E34').Value = '  +0.20%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('B36').Value = 'OKB'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '56.50'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.57%  '
$ws.Range('B37').Value = 'Maker'
$ws.Range('C37').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D37').Value = '3.670.13'
$ws.Range('E37').Value = '  -7.21%  '
$ws.Range('E38').Value = '  -7.22%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '34.19'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +5.05%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.130'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.85%  '
$ws.Range('E41').Value = '  -2.61%  '
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.10'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.55%  '
$ws.Range('B43').Value = 'ApeXProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.35'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.09%  '
$ws.Range('E44').Value = '  -1.64%  '
$ws.Range('D45').Value = '0.0₃0659'
$ws.Range('E45').Value = '  -3.62%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0405'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.97%  '
$ws.Range('E47').Value = '  +1.77%  '
$ws.Range('E48').Value = '  -0.81%  '
$ws.Range('E49').Value = '  -0.27%  '
$ws.Range('E50').Value = '  -1.45%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '129.53'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.89%  '
